$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update footer "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 21:59"

# --- Update country statistics rows ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3398258
$ws.Range("C4").Value = 42612
$ws.Range("D4").Value = 1506863
$ws.Range("E4").Value = 1753715
$ws.Range("G4").Value = 278
$ws.Range("H4").Value = 137680

# Row 6: India
$ws.Range("B6").Value = 879466
$ws.Range("C6").Value = 29108
$ws.Range("D6").Value = 554429
$ws.Range("E6").Value = 301850

# Row 8: Peru
$ws.Range("B8").Value = 326326
$ws.Range("C8").Value = 3616
$ws.Range("D8").Value = 217111
$ws.Range("E8").Value = 97345
$ws.Range("G8").Value = 188
$ws.Range("H8").Value = 11870

# Row 13: Sudafrica
$ws.Range("B13").Value = 276242
$ws.Range("C13").Value = 12058
$ws.Range("D13").Value = 134874
$ws.Range("E13").Value = 137289
$ws.Range("G13").Value = 108
$ws.Range("H13").Value = 4079

# Row 19: Alemania
$ws.Range("B19").Value = 199949
$ws.Range("C19").Value = 137
$ws.Range("E19").Value = 6215

# Row 31: Ecuador
$ws.Range("B31").Value = 67870
$ws.Range("C31").Value = 661
$ws.Range("D31").Value = 30283
$ws.Range("E31").Value = 32540
$ws.Range("G31").Value = 16
$ws.Range("H31").Value = 5047

# Rows 69/70: Uzbekistan overtakes Dinamarca in case totals, so the rows swap.
# Row 69 becomes Uzbekistan with updated stats.
$ws.Range("A69").Value = "Uzbekistan"
$ws.Range("B69").Value = 12997
$ws.Range("C69").Value = 484
$ws.Range("D69").Value = 7852
$ws.Range("E69").Value = 5085
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 60

# Row 70 becomes Dinamarca with the (unchanged) stats formerly held by row 69.
$ws.Range("A70").Value = "Dinamarca"
$ws.Range("B70").Value = 12946
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 12077
$ws.Range("E70").Value = 260
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 609
